$wb = $excel.ActiveWorkbook

# Sheets (tab order): 1=tables, 2=admin_menu, 3=admin_functions, 4=admin_roles
$wsTables = $wb.Worksheets.Item(1)
$wsMenu = $wb.Worksheets.Item(2)
$wsFunctions = $wb.Worksheets.Item(3)
$wsRoles = $wb.Worksheets.Item(4)

# admin_menu: flip "status" (column N) from 1 to 0 for rows 3,4,5,7,8
$wsMenu.Range("N3").Value = 0
$wsMenu.Range("N4").Value = 0
$wsMenu.Range("N5").Value = 0
$wsMenu.Range("N7").Value = 0
$wsMenu.Range("N8").Value = 0

# admin_roles: widen column C (was 30.28515625 -> 48.5703125)
$wsRoles.Range("C:C").ColumnWidth = 47.736979166666664

# Update selections: admin_functions gets a new selection at B6
$wsFunctions.Range("B6").Select()

# admin_menu becomes the active/selected tab, with its selection reset to A1
# (this also clears the old topLeftCell/H8 selection state and drops
# "tabSelected" from whichever sheet was active before)
$wsMenu.Range("A1").Select()
$wsMenu.Activate()
